## Fruta / hortaliza, semanal
## Weekly data refresh: a new price observation is inserted at row 239
## (the top of this store/product's date-ordered block), pushing the
## existing rows 239:362 down by one (-> 240:363). The duplicated row
## carries over every other column unchanged; only its date (column D)
## is updated to the new week's date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 239 (copy + insert-shift) so every column besides the
# date keeps the same value as the row it displaced, then stamp the new
# week's date (2022-09-09 -> Excel serial 44813) onto the fresh row.
$ws.Rows.Item(239).Copy()
$ws.Rows.Item(239).Insert()
$ws.Range("D239").Value = 44813
